$d = $word.ActiveDocument

# 1) Replace the old multi-run title text with the new title text.
#    Word's Find/Replace collapses the matched range into a single run that
#    carries the formatting of the start of the match (the existing bold
#    LMRoman12-Bold 20.5pt title formatting), which also removes the extra
#    runs that made up the old title ("interpretability", "of black-box
#    models", "when", "coupled with local explanations").
$old = "Efficacy of the radial tour and application to extend interpretability of black-box models when coupled with local explanations"
$new = "The use of radial tour to understand variable importance, with application to interpreting black-box machine learning models"
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# 2) Remove the now-empty paragraph that used to follow the title.
#    Deleting the title paragraph's own end-of-paragraph mark merges the
#    title paragraph with the following (empty) paragraph; the surviving
#    paragraph mark comes from the empty paragraph, which carries no
#    special run formatting (so the stray <w:rPr> bold/font run properties
#    that were on the title's paragraph mark are dropped too).
$titlePara = $d.Paragraphs.Item(1)
$markRange = $d.Range($titlePara.Range.End - 1, $titlePara.Range.End)
$markRange.Delete()
